$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2
$ws.Range("E2").Value = 0.000293884
$ws.Range("F2").Value = 0.020380019
$ws.Range("G2").Value = 0.00045500844910000005

# Row 3
$ws.Range("E3").Value = 0.004369107
$ws.Range("F3").Value = 0.009086218
$ws.Range("G3").Value = 0.004884637435972629

# Row 4
$ws.Range("E4").Value = 0.006697469
$ws.Range("F4").Value = 0.01555142
$ws.Range("G4").Value = 0.00781831336973479
